$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: replace [code] with tde
$meta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/TypeCarteVS"

# Version: 0.1.0 -> 2.0.0
$meta.Range("B3").Value = "2.0.0"

# Date: update timestamp
$meta.Range("B8").Value = "2026-01-15T15:25:18+00:00"

# --- Include #0 sheet ---
$inc = $wb.Worksheets.Item("Include #0")

# System URI: replace [code] with tde
$inc.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/CodeSystem/type-carte-code-system"
